$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these D-column cells to remain Text so numeric-looking
# strings (e.g. "501.00", "0.140") are not coerced into numbers,
# preserving the original inline-string formatting.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values from the source diff.
$ws.Range("D2").Value = "66.376.60"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "3.215.21"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "608.13"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "156.26"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.213.67"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").Value = "0.547"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "0.501"
$ws.Range("E12").Value = "  -3.60%  "
$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "38.37"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").Value = "3.742.94"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "66.446.20"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "3.213.92"
$ws.Range("D19").Value = "0.113"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "507.48"
$ws.Range("E20").Value = "  -3.34%  "
$ws.Range("D21").Value = "15.27"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").Value = "0.729"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "8.03"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").Value = "14.58"
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("D25").Value = "85.02"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "0.140"
$ws.Range("E27").Value = "  +54.01%  "
$ws.Range("D28").Value = "3.01"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("E29").Value = "  -3.26%  "
$ws.Range("D30").Value = "2.35"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "6.94"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").Value = "28.26"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -5.42%  "
$ws.Range("D36").Value = "6.42"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "55.44"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "501.00"
$ws.Range("E38").Value = "  -3.23%  "
$ws.Range("D39").Value = "0.0₃0774"
$ws.Range("E39").Value = "  +13.27%  "
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("D42").Value = "3.04"
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("D43").Value = "8.73"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D45").Value = "2.924.28"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "2.44"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("D47").Value = "28.10"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D51").Value = "121.74"
$ws.Range("E51").Value = "  +0.09%  "

Write-Host "Updated $($wb.Name): applied latest crypto price/volume snapshot."
